$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 9 (item id 5487)
$ws.Range("H9").Value = 81.28570999999999
$ws.Range("I9").Value = 82.25
$ws.Range("K9").Value = 82.25
$ws.Range("M9").Value = 86.75
# row 11 (item id 5533)
$ws.Range("H11").Value = 22
$ws.Range("I11").Value = 22
$ws.Range("K11").Value = 22
$ws.Range("M11").Value = 118
# row 38 (item id 4599)
$ws.Range("H38").Value = 1776.6154
$ws.Range("I38").Value = 534.5
$ws.Range("J38").Value = 2841.2856
$ws.Range("K38").Value = 1603.5
$ws.Range("L38").Value = 8523.856800000001
$ws.Range("M38").Value = -1231.5
$ws.Range("N38").Value = -9267.856800000001
# row 132 (item id 44049)
$ws.Range("H132").Value = 52324.25
$ws.Range("I132").Value = 52324.25
$ws.Range("K132").Value = 156972.75
$ws.Range("M132").Value = -154442.75
# row 135 (item id 44047)
$ws.Range("H135").Value = 1783.7273
$ws.Range("I135").Value = 1439.4445
$ws.Range("K135").Value = 12955.0005
$ws.Range("M135").Value = -10420.0005
# row 137 (item id 44013)
$ws.Range("H137").Value = 3042.182
$ws.Range("I137").Value = 1200.3334
$ws.Range("J137").Value = 3732.875
$ws.Range("K137").Value = 3601.0002
$ws.Range("L137").Value = 11198.625
$ws.Range("M137").Value = -1051.0002
$ws.Range("N137").Value = -16298.625

$ws = $wb.Worksheets.Item("ARM")
# row 2 (item id 27713)
$ws.Range("H2").Value = 4015.8667
$ws.Range("I2").Value = 4064
$ws.Range("J2").Value = 3943.6667
$ws.Range("K2").Value = 4064
$ws.Range("L2").Value = 3943.6667
$ws.Range("M2").Value = -3951
$ws.Range("N2").Value = -4169.6667
# row 61 (item id 43999)
$ws.Range("H61").Value = 2375
$ws.Range("I61").Value = 2375
$ws.Range("K61").Value = 2375
$ws.Range("M61").Value = -2163
# row 74 (item id 44000)
$ws.Range("H74").Value = 7596.6924
$ws.Range("I74").Value = 7506.45
$ws.Range("K74").Value = 7506.45
$ws.Range("M74").Value = -6632.45
# row 77 (item id 44000)
$ws.Range("H77").Value = 7596.6924
$ws.Range("I77").Value = 7506.45
$ws.Range("K77").Value = 37532.25
$ws.Range("M77").Value = -33164.25
# row 97 (item id 19941)
$ws.Range("H97").Value = 944.6
$ws.Range("I97").Value = 691.6429000000001
$ws.Range("K97").Value = 691.6429000000001
$ws.Range("M97").Value = -195.6429000000001
# row 116 (item id 27713)
$ws.Range("H116").Value = 4015.8667
$ws.Range("I116").Value = 4064
$ws.Range("J116").Value = 3943.6667
$ws.Range("K116").Value = 4064
$ws.Range("L116").Value = 3943.6667
$ws.Range("M116").Value = -1770
$ws.Range("N116").Value = -8531.6667
# row 136 (item id 43999)
$ws.Range("H136").Value = 2375
$ws.Range("I136").Value = 2375
$ws.Range("K136").Value = 7125
$ws.Range("M136").Value = -4575

$ws = $wb.Worksheets.Item("BSM")
# row 3 (item id 27713)
$ws.Range("H3").Value = 4015.8667
$ws.Range("I3").Value = 4064
$ws.Range("J3").Value = 3943.6667
$ws.Range("K3").Value = 4064
$ws.Range("L3").Value = 3943.6667
$ws.Range("M3").Value = -3950
$ws.Range("N3").Value = -4171.6667
# row 20 (item id 14149)
$ws.Range("H20").Value = 1611.6875
$ws.Range("I20").Value = 1659.1333
$ws.Range("K20").Value = 1659.1333
$ws.Range("M20").Value = -1412.1333
# row 92 (item id 18033)
$ws.Range("H92").Value = 39000
$ws.Range("J92").Value = 39000
$ws.Range("L92").Value = 39000
$ws.Range("N92").Value = -43992

$ws = $wb.Worksheets.Item("CRP")
# row 25 (item id 1895)
$ws.Range("H25").Value = 1533.1666
$ws.Range("I25").Value = 1533.1666
$ws.Range("K25").Value = 1533.1666
$ws.Range("M25").Value = -1359.1666
# row 39 (item id 1915)
$ws.Range("H39").Value = 1150
$ws.Range("I39").Value = 1475
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 1475
$ws.Range("L39").Value = 500
$ws.Range("M39").Value = -1084
$ws.Range("N39").Value = -1282
# row 49 (item id 1915)
$ws.Range("H49").Value = 1150
$ws.Range("I49").Value = 1475
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 1475
$ws.Range("L49").Value = 500
$ws.Range("M49").Value = -1293
$ws.Range("N49").Value = -864
# row 56 (item id 1867)
$ws.Range("H56").Value = 193
$ws.Range("I56").Value = 193
$ws.Range("K56").Value = 193
$ws.Range("M56").Value = 652
# row 58 (item id 44021)
$ws.Range("H58").Value = 4352.8
$ws.Range("I58").Value = 2241.3333
$ws.Range("K58").Value = 2241.3333
$ws.Range("M58").Value = -2038.3333
# row 107 (item id 27689)
$ws.Range("H107").Value = 1415.2858
$ws.Range("I107").Value = 1286.3334
$ws.Range("J107").Value = 1512
$ws.Range("K107").Value = 1286.3334
$ws.Range("L107").Value = 1512
$ws.Range("M107").Value = 633.6666
$ws.Range("N107").Value = -5352
# row 134 (item id 44020)
$ws.Range("H134").Value = 2791.7778
$ws.Range("I134").Value = 2814
$ws.Range("K134").Value = 8442
$ws.Range("M134").Value = -5907
# row 136 (item id 44021)
$ws.Range("H136").Value = 4352.8
$ws.Range("I136").Value = 2241.3333
$ws.Range("K136").Value = 6723.999899999999
$ws.Range("M136").Value = -4173.999899999999

$ws = $wb.Worksheets.Item("CUL")
# row 40 (item id 4827)
$ws.Range("H40").Value = 31
$ws.Range("I40").Value = 31.6
$ws.Range("K40").Value = 126.4
$ws.Range("M40").Value = -57.40000000000001
# row 68 (item id 12895)
$ws.Range("H68").Value = 997.25
$ws.Range("J68").Value = 997.25
$ws.Range("L68").Value = 2991.75
$ws.Range("N68").Value = -4613.75
# row 71 (item id 12895)
$ws.Range("H71").Value = 997.25
$ws.Range("J71").Value = 997.25
$ws.Range("L71").Value = 8975.25
$ws.Range("N71").Value = -17087.25
# row 81 (item id 12843)
$ws.Range("H81").Value = 3150
$ws.Range("J81").Value = 3150
$ws.Range("L81").Value = 9450
$ws.Range("N81").Value = -11696
# row 84 (item id 12843)
$ws.Range("H84").Value = 3150
$ws.Range("J84").Value = 3150
$ws.Range("L84").Value = 28350
$ws.Range("N84").Value = -39582

$ws = $wb.Worksheets.Item("GSM")
# row 97 (item id 19940)
$ws.Range("H97").Value = 1200.6666
$ws.Range("I97").Value = 995
$ws.Range("J97").Value = 1241.8
$ws.Range("K97").Value = 995
$ws.Range("L97").Value = 1241.8
$ws.Range("M97").Value = -499
$ws.Range("N97").Value = -2233.8
# row 102 (item id 36169)
$ws.Range("H102").Value = 1606.6666
$ws.Range("I102").Value = 1792.5
$ws.Range("J102").Value = 120
$ws.Range("K102").Value = 1792.5
$ws.Range("L102").Value = 120
$ws.Range("M102").Value = -170.5
$ws.Range("N102").Value = -3364
# row 113 (item id 27710)
$ws.Range("H113").Value = 5590.05
$ws.Range("I113").Value = 3215.2307
$ws.Range("J113").Value = 10000.429
$ws.Range("K113").Value = 3215.2307
$ws.Range("L113").Value = 10000.429
$ws.Range("M113").Value = -1045.2307
$ws.Range("N113").Value = -14340.429
# row 124 (item id 34247)
$ws.Range("H124").Value = 55000
$ws.Range("J124").Value = 55000
$ws.Range("L124").Value = 55000
$ws.Range("N124").Value = -64820

$ws = $wb.Worksheets.Item("LTW")
# row 22 (item id 5277)
$ws.Range("H22").Value = 911.4545000000001
$ws.Range("I22").Value = 818
$ws.Range("J22").Value = 1075
$ws.Range("K22").Value = 818
$ws.Range("L22").Value = 1075
$ws.Range("M22").Value = -523
$ws.Range("N22").Value = -1665
# row 27 (item id 5277)
$ws.Range("H27").Value = 911.4545000000001
$ws.Range("I27").Value = 818
$ws.Range("J27").Value = 1075
$ws.Range("K27").Value = 818
$ws.Range("L27").Value = 1075
$ws.Range("M27").Value = -711
$ws.Range("N27").Value = -1289
# row 68 (item id 12563)
$ws.Range("H68").Value = 7735.1
$ws.Range("I68").Value = 5087.75
$ws.Range("J68").Value = 9500
$ws.Range("K68").Value = 5087.75
$ws.Range("L68").Value = 9500
$ws.Range("M68").Value = -4338.75
$ws.Range("N68").Value = -10998
# row 71 (item id 12563)
$ws.Range("H71").Value = 7735.1
$ws.Range("I71").Value = 5087.75
$ws.Range("J71").Value = 9500
$ws.Range("K71").Value = 25438.75
$ws.Range("L71").Value = 47500
$ws.Range("M71").Value = -21694.75
$ws.Range("N71").Value = -54988

$ws = $wb.Worksheets.Item("WVR")
# row 54 (item id 3413)
$ws.Range("H54").Value = 32795.1
$ws.Range("I54").Value = 4714
$ws.Range("J54").Value = 60876.2
$ws.Range("K54").Value = 4714
$ws.Range("L54").Value = 60876.2
$ws.Range("M54").Value = -4194
$ws.Range("N54").Value = -61916.2
# row 96 (item id 19977)
$ws.Range("H96").Value = 1183.4117
$ws.Range("I96").Value = 1183.4445
$ws.Range("J96").Value = 1183.375
$ws.Range("K96").Value = 1183.4445
$ws.Range("L96").Value = 1183.375
$ws.Range("M96").Value = 189.5554999999999
$ws.Range("N96").Value = -3929.375
# row 122 (item id 36208)
$ws.Range("H122").Value = 904.36365
$ws.Range("I122").Value = 928.381
$ws.Range("K122").Value = 2785.143
$ws.Range("M122").Value = -335.143
